$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append two new rows (6 and 7), mirroring the existing Reverify/Terverifikasi
#     pair pattern (rows 4/5) but for a "Belum Terverifikasi" verification pass ---

# Row 6 mirrors row 4's layout (A=nomorlama, B=nomorbaru, C=Email_Lama, D=Email_Baru hyperlink, ...)
$ws.Range("A4:H4").Copy($ws.Range("A6:H6"))
# Row 7 mirrors row 5's layout (A=nomorbaru, B=nomorlama, C=Email_Baru hyperlink, D=Email_Lama, ...)
$ws.Range("A5:H5").Copy($ws.Range("A7:H7"))

# New status value for both new rows
$ws.Range("G6").Value = "Belum Terverifikasi"
$ws.Range("G7").Value = "Belum Terverifikasi"

# Re-create the mailto hyperlinks on the new rows' email cells (D6, C7), then
# restore their original cell formatting afterwards since Hyperlinks.Add()
# stamps its own style on the target cell.
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:fityanardi24@gmail.com")
$ws.Range("D4").Copy($ws.Range("D6"))

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:fityanardi24@gmail.com")
$ws.Range("C5").Copy($ws.Range("C7"))

# Update the saved selection to reflect the new bottom-most data row
$ws.Range("G8").Select() | Out-Null

# Record the updated absolute path folder for this workbook
$wb.AbsPath = "G:\magangg\ACC-ACCPartner\"
